$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Target cluster" labels in column D (rows 2-5) to reflect the new
# shared-string table ordering introduced by this commit (an "ECs" cluster was
# inserted and "Resolving-Mac" was dropped/renamed).
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("D5").Value = "Neutrophils"

# Row 2 (ECs) updated TPM-derived values
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.165261
$ws.Range("N2").Value = 0.330522
$ws.Range("O2").Value = 0.2628738644684469
$ws.Range("P2").Value = 0.2081763454351005
$ws.Range("Q2").Value = 0.025146940065
$ws.Range("R2").Value = 0.15088164039
$ws.Range("S2").Value = 0.2628738644684469
$ws.Range("T2").Value = 0.2081763454351005

# Row 3 (Inflammatory-Mac) updated TPM-derived values
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1873846666666667
$ws.Range("N3").Value = 0.562154
$ws.Range("O3").Value = 0.2980650696098803
$ws.Range("P3").Value = 0.3540677028812712
$ws.Range("Q3").Value = 0.02851338780333333
$ws.Range("R3").Value = 0.25662049023
$ws.Range("S3").Value = 0.2980650696098803
$ws.Range("T3").Value = 0.3540677028812712

# Row 4 (MuSCs) updated TPM-derived values
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.133048
$ws.Range("N4").Value = 0.266096
$ws.Range("O4").Value = 0.2116339724423665
$ws.Range("P4").Value = 0.1675982016776448
$ws.Range("Q4").Value = 0.02024524892
$ws.Range("R4").Value = 0.12147149352
$ws.Range("S4").Value = 0.2116339724423665
$ws.Range("T4").Value = 0.1675982016776448

# Row 5 (Neutrophils) updated TPM-derived values
$ws.Range("M5").Value = 0.1429766666666667
$ws.Range("N5").Value = 0.42893
$ws.Range("O5").Value = 0.2274270934793063
$ws.Range("P5").Value = 0.2701577500059835
$ws.Range("Q5").Value = 0.02175604448333333
$ws.Range("R5").Value = 0.19580440035
$ws.Range("S5").Value = 0.2274270934793063
$ws.Range("T5").Value = 0.2701577500059835
